# Generate Report for Handoff
#
# The localization job finished translating and is now ready to be handed
# off again. Update the "Status" (Overview's per-language columns, and the
# per-language sheets' "Status" column) from "In Translation" to
# "Ready for handoff", and bump the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-31 15:16:42"

# --- zh-cn sheet ------------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-31 15:16:38"

# --- de-de sheet ------------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-31 15:16:42"

# --- Column widths ----------------------------------------------------------
# The two status/date columns on each sheet grew wider to fit the new,
# longer text ("Ready for handoff" / refreshed timestamps). ColumnWidth is
# quantized by the host to whole-pixel steps, so we pass a value from the
# middle of the input bracket that resolves to the closest attainable width.
$overview.Columns.Item(5).ColumnWidth = 16.33   # column E
$overview.Columns.Item(6).ColumnWidth = 16.33   # column F
$zhcn.Columns.Item(3).ColumnWidth = 16.33        # column C
$dede.Columns.Item(3).ColumnWidth = 16.33        # column C
